# World Country Income Classification.xlsx — fill in the "Income Group"
# label (column A) next to every country row (column B), matching the
# section header that already exists at the top of each block
# (A2="Low Income", A28/A29.."Lower Middle Income", etc.).
#
# Layout recap (1-based rows):
#   Row 1        : header row (Income Group | Country Name)
#   Row 2        : "Low Income" category row (A2 already set) + first country
#   Rows 3-27    : "Low Income" countries            -> need A3:A27
#   Row 28       : "Lower Middle Income" category row (A28 already set)
#   Rows 29-78   : "Lower Middle Income" countries    -> need A29:A78
#   Row 79       : "Upper Middle Income" category row (A79 already set)
#   Rows 80-132  : "Upper Middle Income" countries    -> need A80:A132
#   Row 133      : "High Income" category row (A133 already set)
#   Rows 134-217 : "High Income" countries            -> need A134:A217
#   Row 218      : "Not classified" category row (A218 already set)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:A27").Value = "Low Income"
$ws.Range("A29:A78").Value = "Lower Middle Income"
$ws.Range("A80:A132").Value = "Upper Middle Income"
$ws.Range("A134:A217").Value = "High Income"

# Leave the view/selection where the original author's session ended up.
$ws.Range("G215:G216").Select()
